$d = $word.ActiveDocument

# The document ends with an (empty) paragraph that only holds the
# "_GoBack" bookmark. We need to:
#   1. insert a brand-new paragraph containing "15.10.2024" right before it
#   2. put the long status update text in front of the bookmark inside
#      that last paragraph
#   3. append a single trailing space (as its own run) after the bookmark

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range

# 1. Insert a new paragraph before the bookmark paragraph and give it the
#    date text.
$lastRange.InsertParagraphBefore()
$datePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$datePara.Range.Text = "15.10.2024"

# Re-fetch the (still) last paragraph - it now only contains the bookmark.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkParaRange = $bookmarkPara.Range

# 2. Insert the trailing space first, right before the paragraph mark, so
#    it becomes its own run placed after the bookmark once the main text
#    is inserted in front of the (still collapsed) bookmark.
$endPos = $bookmarkParaRange.End - 1
$spaceRange = $d.Range($endPos, $endPos)
$spaceRange.InsertAfter(" ")

# 3. Insert the long comment text immediately before the bookmark start.
$bm = $d.Bookmarks.Item("_GoBack")
$bmStartRange = $d.Range($bm.Start, $bm.Start)
$bmStartRange.InsertBefore("Zrobiłem pierwsze kroki związane z interakcjami postaci z rzeczą i pisaniem komentarza. Jak wrócę spróbuję coś zrobić. Pamiętaj by chociaż zrobić do projektu.")
